$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "₹599"
$ws.Range("B2").Value = "₹499"
$ws.Range("B3").Value = "₹599"
$ws.Range("B4").Value = "₹599"
$ws.Range("B5").Value = "₹549"
$ws.Range("B6").Value = "₹499"
$ws.Range("B7").Value = "₹599"
$ws.Range("B8").Value = "₹599"
